# Applies the "feat: improve error messages in highpass_filter and apply_padding
# functions, and free allocated memory in main" journal-entry update:
#  - Adds a new Journal entry (row 5) describing the high-pass filtering work.
#  - Switches the active/selected sheet from Progress to Journal.
#  - Updates the remembered selections on both sheets.

$wb = $excel.ActiveWorkbook

$progress = $wb.Worksheets.Item("Progress")
$journal  = $wb.Worksheets.Item("Journal")

# --- New Journal row (row 5): Task / Date / Notes -------------------------
$taskText = "high pass filtering functionality"

$notesText = @"
logic of the high pass filter was not the same with the low pass filter. It was as below.
1. Daryls code layer
- padding was 50 front and back, not like 60 in the low pass filter
- pass on to MATLAB built in conv() function with padded signal and coefficients(which was imported from fir_51.mat file in the project)
2. Matlab Layer
- the conv function just did a 1d convolution with the signal samples and coeffs.
Logic was simple so was easy and straight forward so simple to recreate.
Although I am curious on why the logic was different in both filterings. According to the full report, both just used fir filtering and didnt mention different paddings and post signal zero padding and pre removal were used. I wonder if this is why some misdetections and 
"@

$journal.Range("A5").Value = $taskText
$journal.Range("B5").Value = 45730
$journal.Range("C5").Value = $notesText
$journal.Rows.Item(5).RowHeight = 189.45

# --- Selection on the Progress sheet moves to B10, and it is no longer the
#     tab shown when the workbook is reopened -------------------------------
$progress.Range("B10").Select()

# --- Journal becomes the active/selected sheet, scrolled & selecting C5 ----
$journal.Activate()
$journal.Range("C5").Select()
